# Generate Report for Handoff
# Adds two new localization-status rows (2b00c18e-...md and 90e65a5c-...md)
# to the Overview / zh-cn / de-de sheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

# Row 6 - 2b00c18e-c46b-445f-bdbb-e291ceefcf79.md
$wsOverview.Range("A6").Value = "2b00c18e-c46b-445f-bdbb-e291ceefcf79.md"
$wsOverview.Range("B6").Value = "e2e\2b00c18e-c46b-445f-bdbb-e291ceefcf79.md"
$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-08-26 12:41:52"
$wsOverview.Range("G6").NumberFormat = $dateFmt
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/2b00c18e-c46b-445f-bdbb-e291ceefcf79.md", "", "", "e2e\2b00c18e-c46b-445f-bdbb-e291ceefcf79.md")
$wsOverview.Range("B6").Font.Color = 15570276

# Row 7 - 90e65a5c-e474-4fac-ae46-5da719d0632e.md
$wsOverview.Range("A7").Value = "90e65a5c-e474-4fac-ae46-5da719d0632e.md"
$wsOverview.Range("B7").Value = "e2e\90e65a5c-e474-4fac-ae46-5da719d0632e.md"
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-08-26 12:41:52"
$wsOverview.Range("G7").NumberFormat = $dateFmt
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/90e65a5c-e474-4fac-ae46-5da719d0632e.md", "", "", "e2e\90e65a5c-e474-4fac-ae46-5da719d0632e.md")
$wsOverview.Range("B7").Font.Color = 15570276

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)

$loZhCn.ListRows.Add() | Out-Null
$loZhCn.ListRows.Add() | Out-Null

# Row 6 - 2b00c18e-c46b-445f-bdbb-e291ceefcf79.md
$wsZhCn.Range("A6").Value = "2b00c18e-c46b-445f-bdbb-e291ceefcf79.md"
$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("F6").Value = "False"
$wsZhCn.Range("G6").Value = "2b00c18e-c46b-445f-bdbb-e291ceefcf79.e175be10d16b92a0c2437770df01995ca849dba4.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2016-08-26 12:41:47"
$wsZhCn.Range("H6").NumberFormat = $dateFmt
$wsZhCn.Range("K6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K6").NumberFormat = $dateFmt
$wsZhCn.Range("M6").Value = "True"
$wsZhCn.Range("O6").Value = "False"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/2b00c18e-c46b-445f-bdbb-e291ceefcf79.md", "", "", "2b00c18e-c46b-445f-bdbb-e291ceefcf79.md")
$wsZhCn.Range("A6").Font.Color = 15570276

# Row 7 - 90e65a5c-e474-4fac-ae46-5da719d0632e.md
$wsZhCn.Range("A7").Value = "90e65a5c-e474-4fac-ae46-5da719d0632e.md"
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "False"
$wsZhCn.Range("G7").Value = "90e65a5c-e474-4fac-ae46-5da719d0632e.2d685bdebe57fb22bf6dbb0319eb11af379f511a.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-08-26 12:41:47"
$wsZhCn.Range("H7").NumberFormat = $dateFmt
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K7").NumberFormat = $dateFmt
$wsZhCn.Range("M7").Value = "True"
$wsZhCn.Range("O7").Value = "False"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/90e65a5c-e474-4fac-ae46-5da719d0632e.md", "", "", "90e65a5c-e474-4fac-ae46-5da719d0632e.md")
$wsZhCn.Range("A7").Font.Color = 15570276

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)

$loDeDe.ListRows.Add() | Out-Null
$loDeDe.ListRows.Add() | Out-Null

# Row 6 - 2b00c18e-c46b-445f-bdbb-e291ceefcf79.md
$wsDeDe.Range("A6").Value = "2b00c18e-c46b-445f-bdbb-e291ceefcf79.md"
$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("F6").Value = "False"
$wsDeDe.Range("G6").Value = "2b00c18e-c46b-445f-bdbb-e291ceefcf79.e175be10d16b92a0c2437770df01995ca849dba4.de-de.xlf"
$wsDeDe.Range("H6").Value = "2016-08-26 12:41:52"
$wsDeDe.Range("H6").NumberFormat = $dateFmt
$wsDeDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K6").NumberFormat = $dateFmt
$wsDeDe.Range("M6").Value = "True"
$wsDeDe.Range("O6").Value = "False"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/2b00c18e-c46b-445f-bdbb-e291ceefcf79.md", "", "", "2b00c18e-c46b-445f-bdbb-e291ceefcf79.md")
$wsDeDe.Range("A6").Font.Color = 15570276

# Row 7 - 90e65a5c-e474-4fac-ae46-5da719d0632e.md
$wsDeDe.Range("A7").Value = "90e65a5c-e474-4fac-ae46-5da719d0632e.md"
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "False"
$wsDeDe.Range("G7").Value = "90e65a5c-e474-4fac-ae46-5da719d0632e.2d685bdebe57fb22bf6dbb0319eb11af379f511a.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-08-26 12:41:52"
$wsDeDe.Range("H7").NumberFormat = $dateFmt
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K7").NumberFormat = $dateFmt
$wsDeDe.Range("M7").Value = "True"
$wsDeDe.Range("O7").Value = "False"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/90e65a5c-e474-4fac-ae46-5da719d0632e.md", "", "", "90e65a5c-e474-4fac-ae46-5da719d0632e.md")
$wsDeDe.Range("A7").Font.Color = 15570276

Write-Output "Added handoff rows for 2b00c18e-...md and 90e65a5c-...md to Overview, zh-cn, de-de sheets."
